$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 86: repurpose the old blank row into a new activity-log entry ---
$b86 = $ws.Range("B86")
$b86.Value = "801 - Comments - Comment modified"
$b86.Font.Name = "Arial"
$b86.Font.Color = 6710886
$b86.Font.Size = 11
$b86.WrapText = $true
$b86.HorizontalAlignment = -4131
$b86.VerticalAlignment = -4160

# --- Row 87: repurpose the old "Save" row into a new activity-log entry ---
$b87 = $ws.Range("B87")
$b87.Value = "802 - Comments - Comment deleted"
$b87.Font.Name = "Arial"
$b87.Font.Color = 6710886
$b87.Font.Size = 11
$b87.WrapText = $true
$b87.HorizontalAlignment = -4131
$b87.VerticalAlignment = -4160

# --- Rows 88-99: twelve new blank filler rows ---
for ($r = 88; $r -le 99; $r++) {
    $c = $ws.Cells.Item($r, 2)
    $c.Font.Name = "Calibri"
    $c.Font.Color = 6710886
    $c.Font.Size = 11
    $c.WrapText = $true
    $c.HorizontalAlignment = -4131
    $c.VerticalAlignment = -4160
}

# --- sheet view: scroll/zoom/selection as left by the editor ---
$ws.Application.ActiveWindow.Zoom = 145
$ws.Application.ActiveWindow.ScrollRow = 70
$ws.Application.ActiveWindow.ScrollColumn = 1
$b86.Select()
